# "added data for 3mL"
#
# Inserts four new worksheets (1000uL, 500uL, 200uL_3mL, 300uL printed data)
# in front of the existing four sheets, matching the target tab order:
#   1000uL_printed.csv, 500uL_printed.csv, 200uL_printed_3mL.csv,
#   300uL_printed.csv, 200uL_printed.csv, 50uL_printed.csv,
#   20uL_printed.csv, 10uL_printed.csv

$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)

# Create in an order such that each Add(firstSheet)/Add(after-previous) call
# produces the exact sheetId sequence and tab order Excel itself assigned:
#   300uL -> id 5 (created first, ends up last of the new block)
#   1000uL -> id 6 (inserted right before 300uL)
#   500uL  -> id 7 (inserted right after 1000uL)
#   200uL_3mL -> id 8 (inserted right after 500uL)
$sNew1 = $wb.Worksheets.Add($firstSheet)
$sNew1.Name = "300uL_printed.csv"

$sNew2 = $wb.Worksheets.Add($sNew1)
$sNew2.Name = "1000uL_printed.csv"

$sNew3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sNew2)
$sNew3.Name = "500uL_printed.csv"

$sNew4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sNew3)
$sNew4.Name = "200uL_printed_3mL.csv"

# NOTE: `Worksheets.Add`'s return value is a *positional* handle (tracks
# "whatever sheet is at that tab index"), not a stable reference to the
# sheet just created -- once a later Add() shifts tab positions, an
# earlier variable silently starts pointing at a different sheet. Re-fetch
# each sheet by its (now-final) name before writing to it so every
# variable below is unambiguous.
$s1000 = $wb.Worksheets.Item("1000uL_printed.csv")
$s500 = $wb.Worksheets.Item("500uL_printed.csv")
$s200_3mL = $wb.Worksheets.Item("200uL_printed_3mL.csv")
$s300 = $wb.Worksheets.Item("300uL_printed.csv")

# --- 1000uL_printed.csv ---
$s1000.Range("A1").Value = "tip1"
$s1000.Range("B1").Value = "tip2"
$s1000.Range("C1").Value = "tip3"
$s1000.Range("A2").Value = 943.8
$s1000.Range("B2").Value = 946.9
$s1000.Range("C2").Value = 947.3
$s1000.Range("A3").Value = 943.5
$s1000.Range("B3").Value = 944.5
$s1000.Range("C3").Value = 947.8
$s1000.Range("A4").Value = 946.6
$s1000.Range("B4").Value = 946.9
$s1000.Range("C4").Value = 944.8
$s1000.Range("A5").Value = 946.9
$s1000.Range("B5").Value = 946.4
$s1000.Range("C5").Value = 947
$s1000.Range("A6").Value = 947.1
$s1000.Range("B6").Value = 947
$s1000.Range("C6").Value = 938.9

# --- 500uL_printed.csv ---
$s500.Range("A1").Value = "tip1"
$s500.Range("B1").Value = "tip2"
$s500.Range("C1").Value = "tip3"
$s500.Range("A2").Value = 475.9
$s500.Range("B2").Value = 478.5
$s500.Range("C2").Value = 469.3
$s500.Range("A3").Value = 473.6
$s500.Range("B3").Value = 478.9
$s500.Range("C3").Value = 470.3
$s500.Range("A4").Value = 474.1
$s500.Range("B4").Value = 478.3
$s500.Range("C4").Value = 470.1
$s500.Range("A5").Value = 473.8
$s500.Range("B5").Value = 477.8
$s500.Range("C5").Value = 469.3
$s500.Range("A6").Value = 474.1
$s500.Range("B6").Value = 480.8
$s500.Range("C6").Value = 468

# --- 200uL_printed_3mL.csv ---
$s200_3mL.Range("A1").Value = "tip1"
$s200_3mL.Range("B1").Value = "tip2"
$s200_3mL.Range("C1").Value = "tip3"
$s200_3mL.Range("A2").Value = 185.3
$s200_3mL.Range("B2").Value = 183.8
$s200_3mL.Range("C2").Value = 187.1
$s200_3mL.Range("A3").Value = 184.4
$s200_3mL.Range("B3").Value = 185.3
$s200_3mL.Range("C3").Value = 187.5
$s200_3mL.Range("A4").Value = 185.7
$s200_3mL.Range("B4").Value = 185.2
$s200_3mL.Range("C4").Value = 188.1
$s200_3mL.Range("A5").Value = 184.3
$s200_3mL.Range("B5").Value = 183.8
$s200_3mL.Range("C5").Value = 187.3
$s200_3mL.Range("A6").Value = 187.3
$s200_3mL.Range("B6").Value = 186
$s200_3mL.Range("C6").Value = 186.6

# --- 300uL_printed.csv ---
$s300.Range("A1").Value = "tip1"
$s300.Range("B1").Value = "tip2"
$s300.Range("C1").Value = "tip3"
$s300.Range("A2").Value = 283.5
$s300.Range("B2").Value = 286.2
$s300.Range("C2").Value = 285.2
$s300.Range("A3").Value = 285.6
$s300.Range("B3").Value = 286.5
$s300.Range("C3").Value = 286.8
$s300.Range("A4").Value = 286
$s300.Range("B4").Value = 285.3
$s300.Range("C4").Value = 285.4
$s300.Range("A5").Value = 285.5
$s300.Range("B5").Value = 287.2
$s300.Range("C5").Value = 286.4
$s300.Range("A6").Value = 286.1
$s300.Range("B6").Value = 284.9
$s300.Range("C6").Value = 286.7

# --- Selections per sheet, matching the target saved state ---
# (wrapped in [void] - Activate()/Select() return booleans that would
# otherwise be echoed to the output stream)
[void]$s200_3mL.Activate()
[void]$s200_3mL.Range("A1:C6").Select()

[void]$s500.Activate()
[void]$s500.Range("A1:C6").Select()

[void]$s300.Activate()
[void]$s300.Range("D1").Select()

# 1000uL_printed.csv ends up the first tab and the active/selected one
[void]$s1000.Activate()
[void]$s1000.Range("A1:C6").Select()

# --- Window chrome (best effort; engine may not persist these) ---
try {
    $win = $excel.ActiveWindow
    $win.Left = 4820
    $win.Top = 0
    $win.Width = 21520
    $win.Height = 12440
    $win.TabRatio = 0.758
} catch {
}
